$wb = $excel.ActiveWorkbook

# Cypher used for the CaseDetailStat query-stat block in the message log
# (same wording as the existing StatOutput query).
$statCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Akita']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$cypherOutputCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN ['Akita'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``
"

# --- Add the two new worksheets at the end of the workbook ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$caseDetailStat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$caseDetailStat.Name = "CaseDetailStat"

$caseDetailStatMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $caseDetailStat)
$caseDetailStatMsg.Name = "CaseDetailStat_Message"

# --- CaseDetailStat: case-level detail table stats (same layout/values
#     as StatOutput) -------------------------------------------------------
$caseDetailStat.Range("A1").Value = "number_of_files"
$caseDetailStat.Range("B1").Value = "number_of_sample"
$caseDetailStat.Range("C1").Value = "number_of_cases"
$caseDetailStat.Range("D1").Value = "number_of_study"

# Force these as text (not numbers) so they store as plain values, then
# strip the number-format back off so no extra cell style lingers.
$caseDetailStat.Range("A2:D2").NumberFormat = "@"
$caseDetailStat.Range("A2").Value = "1"
$caseDetailStat.Range("B2").Value = "2"
$caseDetailStat.Range("C2").Value = "1"
$caseDetailStat.Range("D2").Value = "1"
$caseDetailStat.Range("A2:D2").ClearFormats()

# --- CaseDetailStat_Message: connection/run log, one block per query
#     execution (CypherOutput-style block + two StatOutput-style blocks,
#     the extra one for the new case-detail query) ------------------------
$row = 1
for ($block = 1; $block -le 3; $block++) {
    $caseDetailStatMsg.Cells.Item($row, 1).Value = "Neo4j_URL:"; $row++
    $caseDetailStatMsg.Cells.Item($row, 1).Value = "bolt://ncias-q2251-c.nci.nih.gov:7687"; $row++
    $caseDetailStatMsg.Cells.Item($row, 1).Value = "User_name:"; $row++
    $caseDetailStatMsg.Cells.Item($row, 1).Value = "neo4j"; $row++
    $caseDetailStatMsg.Cells.Item($row, 1).Value = "PWD:"; $row++
    $caseDetailStatMsg.Cells.Item($row, 1).Value = "icdcDBneo4j0"; $row++
    $caseDetailStatMsg.Cells.Item($row, 1).Value = "Cypher:"; $row++
    if ($block -eq 1) {
        $caseDetailStatMsg.Cells.Item($row, 1).Value = $cypherOutputCypher
    } else {
        $caseDetailStatMsg.Cells.Item($row, 1).Value = $statCypher
    }
    $row++
    $caseDetailStatMsg.Cells.Item($row, 1).Value = "Output:"; $row++
    $caseDetailStatMsg.Cells.Item($row, 1).Value = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC01_Canine_Filter_Breed-Akita_Neo4jData.xlsx"; $row++
}

# The multi-line cypher text (block 1) auto-expands its row height; snap
# it back to the sheet's standard row height so no stray row-height
# override is persisted.
$caseDetailStatMsg.Rows.Item(8).AutoFit()

# Creating sheets via COM activates the last one added; restore the
# original active sheet/tab so the workbook-level view state is unchanged.
$wb.Worksheets.Item(1).Activate()
